$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.637.87"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.633.93"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.17"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.29"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("E11").Value = "  -3.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.865.46"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.633.74"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.18"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.625.37"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.33"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.63"
$ws.Range("E22").Value = "  +3.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.36"
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("E24").Value = "  +3.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.83"
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.27"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.477.58"
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.956"
$ws.Range("E37").Value = "  +6.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.879"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.559"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.75"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.20"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.32"
$ws.Range("E46").Value = "  -5.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.775.00"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.77"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0992"
$ws.Range("E51").Value = "  -0.08%  "
